# Auto-update draw results: append the 2025-10-11 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

# Force the numeric-looking columns (date, phase code, timestamp) to be
# stored as text, matching every other row in the sheet (none of the
# "Date"/"Phase"/"InsertedAt" values are real numbers/dates, they're text).
$ws.Range("A" + $row).NumberFormat = "@"
$ws.Range("C" + $row).NumberFormat = "@"
$ws.Range("E" + $row).NumberFormat = "@"

$ws.Range("A" + $row).Value = "2025-10-11"
$ws.Range("B" + $row).Value = "Pick 4"
$ws.Range("C" + $row).Value = "251011"
$ws.Range("D" + $row).Value = "4-4-2-3"
$ws.Range("E" + $row).Value = "2025-10-11T21:35:04.524+04:00"
